{"js": "const body = context.document.body;\n\n// Load all top-level paragraphs so we can find the title paragraph and the\n// trailing duplicated \"Play ...\" / \"Read our review ...\" paragraphs by text.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text\");\n}\nawait context.sync();\n\nconst titleText =\n  \"Play Burning Slots Cash Mesh for Free - Exciting Bonus Features\";\nconst oldMetaText =\n  \"Read our review of Burning Slots Cash Mesh, play now for free and enjoy exciting bonus features including Wild Symbol, Cash Mesh Feature and Hold and Win Bonus Game.\";\nconst newFaqText =\n  'Can I play \"Burning Slots Cash Mesh\" without registration? Yes, you can try out the demo version without registration. What is the maximum payout for this game? The maximum payout is 50,000x your bet. What is the minimum and maximum bet amount? The minimum bet is \u20ac0.10, and the maximum bet amount is \u20ac50. How many paylines does \"Burning Slots Cash Mesh\" have? It has five paylines.';\n\nlet titleParagraph = null;\nlet trailingTitleParagraph = null;\nlet trailingMetaParagraph = null;\nlet seenTitleOnce = false;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text === titleText) {\n    if (!seenTitleOnce) {\n      titleParagraph = p;\n      seenTitleOnce = true;\n    } else {\n      trailingTitleParagraph = p;\n    }\n  }\n  if (p.text === oldMetaText) {\n    trailingMetaParagraph = p;\n  }\n}\n\n// ---------------------------------------------------------------------\n// 1) Insert a new \"Meta description\" paragraph right after the document\n//    title (Heading1) paragraph. A fresh paragraph inherits the Heading1\n//    style from its neighbour, so rewrite it from raw OOXML (flat-OPC)\n//    to land a plain, unstyled paragraph with the bold label + regular\n//    text runs, matching the target markup exactly.\n// ---------------------------------------------------------------------\nconst metaParagraph = titleParagraph.insertParagraph(\"\", \"After\");\nawait context.sync();\n\nconst metaOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>\" +\n  \"<w:r><w:t>: Read our review of Burning Slots Cash Mesh, play now for free and enjoy exciting bonus features including Wild Symbol, Cash Mesh Feature and Hold and Win Bonus Game.</w:t></w:r></w:p></w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData></pkg:part></pkg:package>\";\n\nconst metaRange = metaParagraph.getRange(\"Whole\");\nmetaRange.insertOoxml(metaOoxml, \"Replace\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 2) Delete the trailing, duplicated bold \"Play Burning Slots...\" title\n//    paragraph entirely.\n// ---------------------------------------------------------------------\ntrailingTitleParagraph.delete();\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// 3) Replace the text of the trailing italic paragraph with the new FAQ\n//    text, keeping its italic formatting and leading empty run.\n// ---------------------------------------------------------------------\nconst metaTextRange = trailingMetaParagraph.getRange();\nmetaTextRange.insertText(newFaqText, \"Replace\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------------\n# 1) Insert a new \"Meta description\" paragraph right after the document\n#    title (the first, Heading1 styled, paragraph).\n# ---------------------------------------------------------------------------\n$titlePara = $d.Paragraphs.Item(1)\n$titleRange = $titlePara.Range\n$titleRange.Collapse(0)              # wdCollapseEnd\n$titleRange.InsertParagraphAfter()\n\n$metaPara = $d.Paragraphs.Item(2)\n$metaRange = $metaPara.Range\n\n$metaXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Burning Slots Cash Mesh, play now for free and enjoy exciting bonus features including Wild Symbol, Cash Mesh Feature and Hold and Win Bonus Game.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$metaRange.InsertXML($metaXml)\n\n# ---------------------------------------------------------------------------\n# 2) Near the end of the document, drop the duplicated bold\n#    \"Play Burning Slots Cash Mesh for Free - Exciting Bonus Features\"\n#    paragraph, and replace the following italic paragraph's text with the\n#    new FAQ blurb (keeping its italic run formatting intact).\n# ---------------------------------------------------------------------------\n$titleText = \"Play Burning Slots Cash Mesh for Free - Exciting Bonus Features\"\n$oldMetaText = \"Read our review of Burning Slots Cash Mesh, play now for free and enjoy exciting bonus features including Wild Symbol, Cash Mesh Feature and Hold and Win Bonus Game.\"\n\n$trailingTitleIndex = -1\n$trailingMetaIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n  if ($t -eq $titleText -and $i -gt 1) {\n    $trailingTitleIndex = $i\n  }\n  if ($t -eq $oldMetaText) {\n    $trailingMetaIndex = $i\n  }\n}\n\nif ($trailingTitleIndex -gt 0) {\n  $d.Paragraphs.Item($trailingTitleIndex).Range.Delete()\n}\n\n# Re-resolve the FAQ paragraph's index since the delete above shifts indices.\n$count = $d.Paragraphs.Count\n$trailingMetaIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n  if ($t -eq $oldMetaText) {\n    $trailingMetaIndex = $i\n  }\n}\n\n$faqPara = $d.Paragraphs.Item($trailingMetaIndex)\n$faqRange = $faqPara.Range\n$faqRange.MoveEnd(1, -1)   # wdCharacter; exclude the trailing paragraph mark\n$faqRange.Text = 'Can I play \"Burning Slots Cash Mesh\" without registration? Yes, you can try out the demo version without registration. What is the maximum payout for this game? The maximum payout is 50,000x your bet. What is the minimum and maximum bet amount? The minimum bet is \u20ac0.10, and the maximum bet amount is \u20ac50. How many paylines does \"Burning Slots Cash Mesh\" have? It has five paylines.'\n"}
